$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update move row (row 2): D2 8 -> 2
$ws.Range("D2").Value = 2

# Update carry row (row 3): add C3 = 50, D3 0 -> 2
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = 2

# Update work row (row 4): add C4 = 100, D4 0 -> 1
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 1

# Update attack row (row 5): D5 3 -> 0
$ws.Range("D5").Value = 0

# Update tough row (row 8): D8 5 -> 0
$ws.Range("D8").Value = 0

# Recalculate formulas so dependent cached values update
$excel.Calculate()

# Update selection to C5
$ws.Range("C5").Select()
